# Remove unused light control files:
# - D8 ("LIGHT_PD_10.17.20") is no longer needed, its text is replaced with
#   "OD_LIGHT_SWITCHING"
# - D6 keeps showing "ODLED_only" (text unchanged, but it now shares the
#   string that used to live only at D8 after the old unused string is
#   dropped from the shared strings table)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "ODLED_only"
$ws.Range("D8").Value = "OD_LIGHT_SWITCHING"

$ws.Range("D13").Select()
